# Aptitudo Question Generation and Send Invitation
# Replaces the old "S2/S3" placeholder question rows and leftover
# "wert"/"test-code-004" test rows on Sheet2 with real invitee data
# (Name / Last name / Email) and turns the e-mail cells into mailto
# hyperlinks, resizing things to suit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

function RGBColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# ---------------------------------------------------------------
# 1. Drop the two trailing leftover test rows (old rows 20 & 21:
#    "wert" / "test-code-004"). Deleting the same row index twice
#    removes both, shifting dimension down to row 19.
# ---------------------------------------------------------------
$ws.Rows(20).Delete()
$ws.Rows(20).Delete()

# ---------------------------------------------------------------
# 2. Replace row 16 ("S2" / "What is your approach ...") with the
#    first invitee: Tamima Tarin <tarin.aiub@gmail.com>
# ---------------------------------------------------------------
$ws.Range("A16").Value = "Tamima"
$ws.Range("A16").Font.Name = "Arial"
$ws.Range("A16").Font.Color = (RGBColor 0x11 0x18 0x27)
$ws.Range("A16").VerticalAlignment = -4108

$ws.Range("B16").Value = "Tarin"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").VerticalAlignment = -4108

$ws.Range("C16").Value = "tarin.aiub@gmail.com"
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:tarin.aiub@gmail.com") | Out-Null

$ws.Range("D16").Clear()
$ws.Rows(16).AutoFit()

# ---------------------------------------------------------------
# 3. Replace row 17 ("S3" / "How do you handle critical bugs...")
#    with the second invitee: Nourin  Ahmed <nourinahmed.nuba@gmail.com>
# ---------------------------------------------------------------
$ws.Range("A17").Value = "Nourin "
$ws.Range("A17").Font.Name = "Arial"
$ws.Range("A17").Font.Color = (RGBColor 0x11 0x18 0x27)
$ws.Range("A17").VerticalAlignment = -4108

$ws.Range("B17").Value = "Ahmed"
$ws.Range("B17").WrapText = $true
$ws.Range("B17").VerticalAlignment = -4108

$ws.Range("C17").Value = "nourinahmed.nuba@gmail.com"
$ws.Range("C17").WrapText = $true
$ws.Range("C17").VerticalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:nourinahmed.nuba@gmail.com") | Out-Null

$ws.Range("D17").Clear()
$ws.Rows(17).RowHeight = 30

# ---------------------------------------------------------------
# 4. Rows 18 & 19 become blank placeholder rows (keep their
#    existing Arial/FF111827 style on column A, add a new
#    Arial/FF222222 style on column C).
# ---------------------------------------------------------------
$ws.Range("A18").ClearContents()
$ws.Range("C18").Font.Name = "Arial"
$ws.Range("C18").Font.Color = (RGBColor 0x22 0x22 0x22)

$ws.Range("A19").ClearContents()

# ---------------------------------------------------------------
# 5. Column C needs to be wider to comfortably show e-mail
#    addresses.
# ---------------------------------------------------------------
$ws.Columns(3).ColumnWidth = 26.5

# ---------------------------------------------------------------
# 6. Update the view: scroll/zoom/selection.
# ---------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$win.Zoom = 98
$ws.Rows(15).Select()
